$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H4").Value = 412.5
$ws.Range("J4").Value = 600
$ws.Range("L4").Value = 600
$ws.Range("N4").Value = -828

$ws.Range("H6").Value = 1220.6666
$ws.Range("I6").Value = 274.66666
$ws.Range("J6").Value = 2166.6667
$ws.Range("K6").Value = 823.9999799999999
$ws.Range("L6").Value = 6500.000100000001
$ws.Range("M6").Value = -711.9999799999999
$ws.Range("N6").Value = -6724.000100000001

$ws.Range("H39").Value = 600.8
$ws.Range("I39").Value = 28
$ws.Range("J39").Value = 1102
$ws.Range("K39").Value = 84
$ws.Range("L39").Value = 3306
$ws.Range("M39").Value = 212
$ws.Range("N39").Value = -3898

$ws.Range("H86").Value = 6870.9375
$ws.Range("I86").Value = 8981.333000000001
$ws.Range("K86").Value = 8981.333000000001
$ws.Range("M86").Value = -7858.333000000001

$ws.Range("H89").Value = 6870.9375
$ws.Range("I89").Value = 8981.333000000001
$ws.Range("K89").Value = 44906.665
$ws.Range("M89").Value = -39290.665

$ws.Range("H92").Value = 750
$ws.Range("I92").Value = 642.8570999999999
$ws.Range("K92").Value = 642.8570999999999
$ws.Range("M92").Value = 605.1429000000001

$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -44992

$ws.Range("H99").Value = 1406.3334
$ws.Range("I99").Value = 400.83334
$ws.Range("J99").Value = 2411.8333
$ws.Range("K99").Value = 1202.50002
$ws.Range("L99").Value = 7235.499899999999
$ws.Range("M99").Value = 295.4999800000001
$ws.Range("N99").Value = -10231.4999

$ws.Range("H100").Value = 5187
$ws.Range("I100").Value = 4897.091
$ws.Range("J100").Value = 5541.3335
$ws.Range("K100").Value = 4897.091
$ws.Range("L100").Value = 5541.3335
$ws.Range("M100").Value = -4356.091
$ws.Range("N100").Value = -6623.3335

$ws.Range("H101").Value = 1108.5
$ws.Range("I101").Value = 264
$ws.Range("J101").Value = 1530.75
$ws.Range("K101").Value = 792
$ws.Range("L101").Value = 4592.25
$ws.Range("M101").Value = 830
$ws.Range("N101").Value = -7836.25

$ws.Range("H132").Value = 4652957.5
$ws.Range("I132").Value = 5001647
$ws.Range("K132").Value = 15004941
$ws.Range("M132").Value = -15002411


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H4").Value = 260.42856
$ws.Range("I4").Value = 184.6
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 184.6
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -68.59999999999999
$ws.Range("N4").Value = -682

$ws.Range("H5").Value = 172.44444
$ws.Range("I5").Value = 72.2
$ws.Range("J5").Value = 297.75
$ws.Range("K5").Value = 72.2
$ws.Range("L5").Value = 297.75
$ws.Range("M5").Value = 39.8
$ws.Range("N5").Value = -521.75

$ws.Range("H34").Value = 47280.363
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 47280.363
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 47280.363
$ws.Range("N34").Value = -47822.363
$ws.Range("M34").ClearContents()

$ws.Range("H92").Value = 35540
$ws.Range("J92").Value = 35540
$ws.Range("L92").Value = 35540
$ws.Range("N92").Value = -40532

$ws.Range("H94").Value = 30023.357
$ws.Range("J94").Value = 30023.357
$ws.Range("L94").Value = 30023.357
$ws.Range("N94").Value = -31825.357

$ws.Range("H137").Value = 44000
$ws.Range("J137").Value = 44000
$ws.Range("L137").Value = 44000
$ws.Range("N137").Value = -54200


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 172.44444
$ws.Range("I4").Value = 72.2
$ws.Range("J4").Value = 297.75
$ws.Range("K4").Value = 72.2
$ws.Range("L4").Value = 297.75
$ws.Range("M4").Value = 42.8
$ws.Range("N4").Value = -527.75

$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 853.1739
$ws.Range("I94").Value = 701.15
$ws.Range("K94").Value = 701.15
$ws.Range("M94").Value = -250.15


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 56.555557
$ws.Range("I7").Value = 43.22222
$ws.Range("J7").Value = 69.888885
$ws.Range("K7").Value = 43.22222
$ws.Range("L7").Value = 69.888885
$ws.Range("M7").Value = 69.77778000000001
$ws.Range("N7").Value = -295.888885

$ws.Range("H74").Value = 17839
$ws.Range("J74").Value = 17839
$ws.Range("L74").Value = 17839
$ws.Range("N74").Value = -19587

$ws.Range("H77").Value = 17839
$ws.Range("J77").Value = 17839
$ws.Range("L77").Value = 53517
$ws.Range("N77").Value = -62253

$ws.Range("H115").Value = 37499
$ws.Range("J115").Value = 37499
$ws.Range("L115").Value = 37499
$ws.Range("N115").Value = -39849


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 14509.228
$ws.Range("I4").Value = 27409.182
$ws.Range("J4").Value = 1609.2727
$ws.Range("K4").Value = 82227.546
$ws.Range("L4").Value = 4827.8181
$ws.Range("M4").Value = -82115.546
$ws.Range("N4").Value = -5051.8181

$ws.Range("H33").Value = 276.4
$ws.Range("I33").Value = 340
$ws.Range("J33").Value = 234
$ws.Range("K33").Value = 2040
$ws.Range("L33").Value = 1404
$ws.Range("M33").Value = -1757
$ws.Range("N33").Value = -1970

$ws.Range("H44").Value = 1007.875
$ws.Range("I44").Value = 594.3333
$ws.Range("J44").Value = 1256
$ws.Range("K44").Value = 1782.9999
$ws.Range("L44").Value = 3768
$ws.Range("M44").Value = -1384.9999
$ws.Range("N44").Value = -4564

$ws.Range("H113").Value = 4348624
$ws.Range("I113").Value = 33333534
$ws.Range("J113").Value = 887.8
$ws.Range("K113").Value = 100000602
$ws.Range("L113").Value = 2663.4
$ws.Range("M113").Value = -99998432
$ws.Range("N113").Value = -7003.4

$ws.Range("H131").Value = 1053.7528
$ws.Range("I131").Value = 1476.8334
$ws.Range("J131").Value = 987.8182
$ws.Range("K131").Value = 4430.5002
$ws.Range("L131").Value = 2963.4546
$ws.Range("M131").Value = 609.4997999999996
$ws.Range("N131").Value = -13043.4546

$ws.Range("H132").Value = 1813
$ws.Range("I132").Value = 1180.3
$ws.Range("J132").Value = 2867.5
$ws.Range("K132").Value = 10622.7
$ws.Range("L132").Value = 25807.5
$ws.Range("M132").Value = -8092.699999999999
$ws.Range("N132").Value = -30867.5

$ws.Range("H134").Value = 2047.5454
$ws.Range("I134").Value = 863.75
$ws.Range("J134").Value = 3468.1
$ws.Range("K134").Value = 2591.25
$ws.Range("L134").Value = 10404.3
$ws.Range("M134").Value = 2478.75
$ws.Range("N134").Value = -20544.3


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H2").Value = 57
$ws.Range("I2").Value = 47.666668
$ws.Range("K2").Value = 47.666668
$ws.Range("M2").Value = 65.333332


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 90910350
$ws.Range("I22").Value = 166666990
$ws.Range("J22").Value = 2390
$ws.Range("K22").Value = 166666990
$ws.Range("L22").Value = 2390
$ws.Range("M22").Value = -166666695
$ws.Range("N22").Value = -2980

$ws.Range("H27").Value = 90910350
$ws.Range("I27").Value = 166666990
$ws.Range("J27").Value = 2390
$ws.Range("K27").Value = 166666990
$ws.Range("L27").Value = 2390
$ws.Range("M27").Value = -166666883
$ws.Range("N27").Value = -2604

$ws.Range("H35").Value = 10177

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996

$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws.Range("H96").Value = 2021.625
$ws.Range("I96").Value = 1894.75
$ws.Range("J96").Value = 2148.5
$ws.Range("K96").Value = 1894.75
$ws.Range("L96").Value = 2148.5
$ws.Range("M96").Value = -521.75
$ws.Range("N96").Value = -4894.5

$ws.Range("H122").Value = 557553.9399999999
$ws.Range("I122").Value = 626529.4399999999
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 1879588.32
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -1877138.32
$ws.Range("N122").Value = -22150

$ws.Range("H132").Value = 185522.67
$ws.Range("I132").Value = 257852.4
$ws.Range("J132").Value = 9218.9375
$ws.Range("K132").Value = 773557.2
$ws.Range("L132").Value = 27656.8125
$ws.Range("M132").Value = -771027.2
$ws.Range("N132").Value = -32716.8125

